# New walk forward, added early stop to hyperparameter tuning
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update Row 3 (RF) ---
$ws.Range("B3").Value = 0.822
$ws.Range("C3").Value = 0.902
$ws.Range("D3").Value = 0.653
$ws.Range("E3").Value = 0.906
$ws.Range("F3").Value = 0.898
$ws.Range("G3").Value = 0.102
$ws.Range("H3").Value = 0.319
$ws.Range("I3").Value = 0.239
$ws.Range("J3").Value = 0.971

# --- Update Row 4 (NN) ---
$ws.Range("E4").Value = 0.785
$ws.Range("F4").Value = 0.767
$ws.Range("G4").Value = 0.234
$ws.Range("H4").Value = 0.484
$ws.Range("I4").Value = 0.35
$ws.Range("J4").Value = 0.9370000000000001

# --- Update Row 5 (RNN) ---
$ws.Range("E5").Value = 0.6889999999999999
$ws.Range("F5").Value = 0.676
$ws.Range("G5").Value = 0.339
$ws.Range("H5").Value = 0.582
$ws.Range("I5").Value = 0.449
$ws.Range("J5").Value = 0.872

# --- Add new Row 6 (Ensemble) ---
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "Ensemble"

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.61
$ws.Range("F6").Value = 0.578
$ws.Range("G6").Value = 0.426
$ws.Range("H6").Value = 0.653
$ws.Range("I6").Value = 0.488
$ws.Range("J6").Value = 0.867
